$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing header/data text typos
$ws.Range("A1").Value = "firstaame"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"

# New column D: alert text header + value
$ws.Range("D1").Value = "alerttext"
$ws.Range("D2").Value = "Customer added successfully"

# Update selection to match target workbook state
$ws.Range("E12").Select()
